$d = $word.ActiveDocument

# Builds a minimal Flat-OPC (WordOpenXML) document wrapper around a body
# fragment so it can be fed to Range.InsertXML(). InsertXML replaces only
# the contents of the target Range, leaving sibling runs (including empty
# runs used as formatting placeholders) completely untouched.
function Get-XmlWrapped($innerXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $innerXml + '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

function Escape-Xml($s) {
    return $s.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
}

# Replaces the Nth (1-based) occurrence of $oldText anywhere in the
# document body with $newText, re-using the supplied run properties XML
# (e.g. "<w:rPr><w:b/></w:rPr>") so that character formatting on the run
# is preserved. Only the matched run's own text is rewritten; any
# neighboring empty runs in the same paragraph are left alone.
function Replace-RunText($oldText, $newText, $rPrXml, $occurrence) {
    $count = 0
    foreach ($p in $d.Paragraphs) {
        $full = $p.Range.Text
        $idx = $full.IndexOf($oldText)
        if ($idx -ge 0) {
            $count++
            if ($count -ne $occurrence) {
                continue
            }
            $start = $p.Range.Start + $idx
            $end = $start + $oldText.Length
            $r = $d.Range($start, $end)
            $escaped = Escape-Xml $newText
            $innerXml = '<w:p><w:r>' + $rPrXml + '<w:t>' + $escaped + '</w:t></w:r></w:p>'
            $r.InsertXML((Get-XmlWrapped $innerXml))
            return
        }
    }
}

# 1. Main page title (Heading1)
Replace-RunText "Play Ghostbusters Triple Slime for Free - Slot Game Review" "Play Ghostbusters Triple Slime for Free" "" 1

# 2-5. "What we like" bullet list
Replace-RunText "Elaborate grid design with nostalgic symbols from the Ghostbusters movie." "Elaborate grid design with tribute to Ghostbusters" "" 1
Replace-RunText "Return to player rate is 96.08% providing a fair winning opportunity." "Symbols include iconic characters and items from the movie" "" 1
Replace-RunText "720 paylines provide players with ample chances of winning payouts." "Fair Return to Player rate of 96.08%" "" 1
Replace-RunText "Simple and easy-to-use interface suitable for beginners." "Expanded grid layout with 720 paylines" "" 1

# 6-7. "What we don't like" bullet list
Replace-RunText "May not be as innovative compared to the previous Ghostbusters slot game." "May not be as innovative as previous Ghostbusters slot game" "" 1
Replace-RunText "Not recommended for high-risk players." "Recommended to start with smaller bets for first-time players" "" 1

# 8. Bold meta-title run near the end of the document
Replace-RunText "Play Ghostbusters Triple Slime for Free - Slot Game Review" "Play Ghostbusters Triple Slime for Free" "<w:rPr><w:b/></w:rPr>" 1

# 9. Italic meta-description run
Replace-RunText "Discover everything you need to know about Ghostbusters Triple Slime. Try it for free and enjoy the thrilling gameplay with 720 paylines and a 96.08% RTP rate." "Read our review of Ghostbusters Triple Slime and play it for free. Enjoy the tribute to the classic movie." "<w:rPr><w:i/></w:rPr>" 1

Write-Output "Done"
